# Apply updated currentAveragePrice / Leve price / profit figures
# across the Moogle_Profits sheets (per scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 31317.646
$ws.Range("I33").Value = 36867.855
$ws.Range("K33").Value = 36867.855
$ws.Range("M33").Value = -36638.855
$ws.Range("H39").Value = 481.57144
$ws.Range("I39").Value = 61.833332
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 185.499996
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 110.500004
$ws.Range("N39").Value = -9592
$ws.Range("H43").Value = 2197
$ws.Range("J43").Value = 2197
$ws.Range("L43").Value = 2197
$ws.Range("N43").Value = -2335
$ws.Range("H62").Value = 2529.2
$ws.Range("I62").Value = 1781.4
$ws.Range("J62").Value = 3277
$ws.Range("K62").Value = 1781.4
$ws.Range("L62").Value = 3277
$ws.Range("M62").Value = -1157.4
$ws.Range("N62").Value = -4525
$ws.Range("H65").Value = 2529.2
$ws.Range("I65").Value = 1781.4
$ws.Range("J65").Value = 3277
$ws.Range("K65").Value = 8907
$ws.Range("L65").Value = 16385
$ws.Range("M65").Value = -5787
$ws.Range("N65").Value = -22625
$ws.Range("H107").Value = 1012.2083
$ws.Range("I107").Value = 925.1
$ws.Range("K107").Value = 925.1
$ws.Range("M107").Value = 994.9
$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -100119
$ws.Range("H141").Value = 6447.4116
$ws.Range("I141").Value = 2372.4285
$ws.Range("K141").Value = 7117.2855
$ws.Range("M141").Value = -1937.2855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3058.6316
$ws.Range("I61").Value = 1775.6
$ws.Range("K61").Value = 1775.6
$ws.Range("M61").Value = -1563.6
$ws.Range("H74").Value = 4224.9062
$ws.Range("I74").Value = 707.6957
$ws.Range("K74").Value = 707.6957
$ws.Range("M74").Value = 166.3043
$ws.Range("H77").Value = 4224.9062
$ws.Range("I77").Value = 707.6957
$ws.Range("K77").Value = 3538.4785
$ws.Range("M77").Value = 829.5214999999998
$ws.Range("H102").Value = 2924.65
$ws.Range("I102").Value = 2621.9443
$ws.Range("K102").Value = 2621.9443
$ws.Range("M102").Value = -999.9443000000001
$ws.Range("H122").Value = 5436.2144
$ws.Range("I122").Value = 3048
$ws.Range("K122").Value = 9144
$ws.Range("M122").Value = -6694
$ws.Range("H123").Value = 82500
$ws.Range("J123").Value = 82500
$ws.Range("L123").Value = 82500
$ws.Range("N123").Value = -92300
$ws.Range("H132").Value = 3792.85
$ws.Range("I132").Value = 1658.6154
$ws.Range("K132").Value = 4975.8462
$ws.Range("M132").Value = -2445.8462
$ws.Range("H136").Value = 3058.6316
$ws.Range("I136").Value = 1775.6
$ws.Range("K136").Value = 5326.799999999999
$ws.Range("M136").Value = -2776.799999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 22853.945
$ws.Range("I80").Value = 60405.6
$ws.Range("J80").Value = 8411
$ws.Range("K80").Value = 60405.6
$ws.Range("L80").Value = 8411
$ws.Range("M80").Value = -59407.6
$ws.Range("N80").Value = -10407
$ws.Range("H83").Value = 22853.945
$ws.Range("I83").Value = 60405.6
$ws.Range("J83").Value = 8411
$ws.Range("K83").Value = 302028
$ws.Range("L83").Value = 42055
$ws.Range("M83").Value = -297036
$ws.Range("N83").Value = -52039
$ws.Range("H86").Value = 4763.6875
$ws.Range("I86").Value = 1893.0834
$ws.Range("J86").Value = 13375.5
$ws.Range("K86").Value = 1893.0834
$ws.Range("L86").Value = 13375.5
$ws.Range("M86").Value = -770.0834
$ws.Range("N86").Value = -15621.5
$ws.Range("H89").Value = 4763.6875
$ws.Range("I89").Value = 1893.0834
$ws.Range("J89").Value = 13375.5
$ws.Range("K89").Value = 9465.416999999999
$ws.Range("L89").Value = 66877.5
$ws.Range("M89").Value = -3849.416999999999
$ws.Range("N89").Value = -78109.5
$ws.Range("H99").Value = 2322.4375
$ws.Range("I99").Value = 1897.6154
$ws.Range("K99").Value = 1897.6154
$ws.Range("M99").Value = -399.6153999999999
$ws.Range("H105").Value = 6831.478
$ws.Range("I105").Value = 4818.0557
$ws.Range("J105").Value = 14079.8
$ws.Range("K105").Value = 4818.0557
$ws.Range("L105").Value = 14079.8
$ws.Range("M105").Value = -3071.0557
$ws.Range("N105").Value = -17573.8
$ws.Range("H107").Value = 2630.4285
$ws.Range("I107").Value = 2591.5264
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2591.5264
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -671.5264000000002
$ws.Range("N107").Value = -6840
$ws.Range("H134").Value = 5761
$ws.Range("I134").Value = 1655.625
$ws.Range("K134").Value = 4966.875
$ws.Range("M134").Value = -2431.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5935
$ws.Range("J58").Value = 9871
$ws.Range("L58").Value = 9871
$ws.Range("N58").Value = -10277
$ws.Range("H105").Value = 2209.4856
$ws.Range("I105").Value = 1997.1305
$ws.Range("K105").Value = 1997.1305
$ws.Range("M105").Value = -250.1305
$ws.Range("H124").Value = 72999
$ws.Range("J124").Value = 72999
$ws.Range("L124").Value = 72999
$ws.Range("N124").Value = -77909
$ws.Range("H132").Value = 2906.9736
$ws.Range("I132").Value = 2424.743
$ws.Range("K132").Value = 7274.228999999999
$ws.Range("M132").Value = -4744.228999999999
$ws.Range("H134").Value = 3511.7778
$ws.Range("I134").Value = 2574.96
$ws.Range("J134").Value = 5640.909
$ws.Range("K134").Value = 7724.88
$ws.Range("L134").Value = 16922.727
$ws.Range("M134").Value = -5189.88
$ws.Range("N134").Value = -21992.727
$ws.Range("H136").Value = 5935
$ws.Range("J136").Value = 9871
$ws.Range("L136").Value = 29613
$ws.Range("N136").Value = -34713

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 9806024
$ws.Range("I129").Value = 619.8333
$ws.Range("J129").Value = 15154427
$ws.Range("K129").Value = 1859.4999
$ws.Range("L129").Value = 45463281
$ws.Range("M129").Value = 3140.5001
$ws.Range("N129").Value = -45473281

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 30458
$ws.Range("I57").Value = 31249.75
$ws.Range("J57").Value = 28874.5
$ws.Range("K57").Value = 31249.75
$ws.Range("L57").Value = 28874.5
$ws.Range("M57").Value = -30429.75
$ws.Range("N57").Value = -30514.5
$ws.Range("H70").Value = 4687.684
$ws.Range("I70").Value = 4612.2856
$ws.Range("J70").Value = 4898.8
$ws.Range("K70").Value = 4612.2856
$ws.Range("L70").Value = 4898.8
$ws.Range("M70").Value = -4342.2856
$ws.Range("N70").Value = -5438.8
$ws.Range("H73").Value = 4687.684
$ws.Range("I73").Value = 4612.2856
$ws.Range("J73").Value = 4898.8
$ws.Range("K73").Value = 4612.2856
$ws.Range("L73").Value = 4898.8
$ws.Range("M73").Value = -3676.2856
$ws.Range("N73").Value = -6770.8
$ws.Range("H80").Value = 5312.75
$ws.Range("J80").Value = 5799.4
$ws.Range("L80").Value = 5799.4
$ws.Range("N80").Value = -7795.4
$ws.Range("H83").Value = 5312.75
$ws.Range("J83").Value = 5799.4
$ws.Range("L83").Value = 28997
$ws.Range("N83").Value = -38981
$ws.Range("H122").Value = 3157.7368
$ws.Range("I122").Value = 1408.7693
$ws.Range("J122").Value = 6947.1665
$ws.Range("K122").Value = 4226.3079
$ws.Range("L122").Value = 20841.4995
$ws.Range("M122").Value = -1776.3079
$ws.Range("N122").Value = -25741.4995
$ws.Range("H126").Value = 4177.923
$ws.Range("I126").Value = 3573.8572
$ws.Range("J126").Value = 4882.6665
$ws.Range("K126").Value = 10721.5716
$ws.Range("L126").Value = 14647.9995
$ws.Range("M126").Value = -8251.571599999999
$ws.Range("N126").Value = -19587.9995
$ws.Range("H132").Value = 3298.8215
$ws.Range("I132").Value = 2754.3333
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 8262.999899999999
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -5732.999899999999
$ws.Range("N132").Value = -59060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1611.8334
$ws.Range("I16").Value = 846.13635
$ws.Range("J16").Value = 10034.5
$ws.Range("K16").Value = 846.13635
$ws.Range("L16").Value = 10034.5
$ws.Range("M16").Value = -676.13635
$ws.Range("N16").Value = -10374.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1452.3103
$ws.Range("I107").Value = 1353.15
$ws.Range("J107").Value = 1672.6666
$ws.Range("K107").Value = 4059.45
$ws.Range("L107").Value = 5017.9998
$ws.Range("M107").Value = -2139.45
$ws.Range("N107").Value = -8857.9998
$ws.Range("H132").Value = 8458.5
$ws.Range("I132").Value = 3749.5
$ws.Range("J132").Value = 13167.5
$ws.Range("K132").Value = 11248.5
$ws.Range("L132").Value = 39502.5
$ws.Range("M132").Value = -8718.5
$ws.Range("N132").Value = -44562.5
$ws.Range("H135").Value = 85357.5
$ws.Range("J135").Value = 85357.5
$ws.Range("L135").Value = 85357.5
$ws.Range("N135").Value = -95497.5

